$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("ECs", "F2", "F2r", "ECs", 3, 1, 0.5597896666666667, 1.679369, 0.4918833646099291, 0.491883364609929, 2, 0.6666666666666666, 8.038446, 24.115338, 0.1223314651384763, 0.1223314651384763, 4.499839006858001, 40.498551061722, 0.06017281266997597, 0.06017281266997596)
    ,("ECs", "F2", "F2r", "FAPs", 3, 1, 0.5597896666666667, 1.679369, 0.4918833646099291, 0.491883364609929, 3, 1, 44.50790266666667, 133.523708, 0.6773345175739228, 0.6773345175739228, 24.91506399780578, 224.235575980252, 0.3331695814707043, 0.3331695814707042)
    ,("ECs", "F2", "F2r", "sCs", 3, 1, 0.5597896666666667, 1.679369, 0.4918833646099291, 0.491883364609929, 3, 1, 13.16402266666667, 39.492068, 0.2003340172876008, 0.2003340172876008, 7.369083860565779, 66.32175474509201, 0.09854097046924877, 0.09854097046924876)
    ,("FAPs", "F2", "F2r", "ECs", 2, 0.6666666666666666, 0.4205383333333333, 1.261615, 0.3695241671379879, 0.3695241671379879, 2, 0.6666666666666666, 8.038446, 24.115338, 0.1223314651384763, 0.1223314651384763, 3.38047468343, 30.42427215087, 0.04520443277006525, 0.04520443277006525)
    ,("FAPs", "F2", "F2r", "FAPs", 2, 0.6666666666666666, 0.4205383333333333, 1.261615, 0.3695241671379879, 0.3695241671379879, 3, 1, 44.50790266666667, 133.523708, 0.6773345175739228, 0.6773345175739228, 18.71727920760222, 168.45551286842, 0.2502914734803146, 0.2502914734803146)
    ,("FAPs", "F2", "F2r", "sCs", 2, 0.6666666666666666, 0.4205383333333333, 1.261615, 0.3695241671379879, 0.3695241671379879, 3, 1, 13.16402266666667, 39.492068, 0.2003340172876008, 0.2003340172876008, 5.535976152202222, 49.82378536982, 0.07402826088760794, 0.07402826088760794)
    ,("sCs", "F2", "F2r", "ECs", 2, 0.6666666666666666, 0.1577256666666667, 0.473177, 0.138592468252083, 0.138592468252083, 2, 0.6666666666666666, 8.038446, 24.115338, 0.1223314651384763, 0.1223314651384763, 1.267869254314, 11.410823288826, 0.01695421969843508, 0.01695421969843507)
    ,("sCs", "F2", "F2r", "FAPs", 2, 0.6666666666666666, 0.1577256666666667, 0.473177, 0.138592468252083, 0.138592468252083, 3, 1, 44.50790266666667, 133.523708, 0.6773345175739228, 0.6773345175739228, 7.020038620035112, 63.180347580316, 0.09387346262290387, 0.09387346262290384)
    ,("sCs", "F2", "F2r", "sCs", 2, 0.6666666666666666, 0.1577256666666667, 0.473177, 0.138592468252083, 0.138592468252083, 3, 1, 13.16402266666667, 39.492068, 0.2003340172876008, 0.2003340172876008, 2.076304251115111, 18.686738260036, 0.02776478593074406, 0.02776478593074405)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
